# Add "Last Location" column (S) to the Device Report template.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Comments: bump the jxls lastCell reference from R7 to S7 -------------
$ws.Range("A1").Comment.Text('jx:area(lastCell="S7")')
$ws.Range("A7").Comment.Text('jx:each(items="devices", var="device", lastCell="S7")')

# --- New column width (matches column R's header/body formatting) ---------
$ws.Columns("S").ColumnWidth = 35.33

# --- Header cell S6: "Last Location" (styled like the other headers) ------
$ws.Range("R6").Copy()
$ws.Range("S6").PasteSpecial(-4122)
$ws.Range("S6").Value = "Last Location"

# --- Data cell S7: the jxls placeholder (styled like the other data cols) -
$ws.Range("R7").Copy()
$ws.Range("S7").PasteSpecial(-4122)
$ws.Range("S7").Value = '${device.lastLocation}'

# --- Selection / scroll position, matching the authored view --------------
$ws.Range("S7").Select()
$excel.ActiveWindow.ScrollColumn = 14
